# Implementing Suggestion given in presentation
#
# Slide 3 ("Our Research Question is ...") - rewrite the quoted research
# question paragraph (5th paragraph of the subtitle placeholder).
# Slide 4 ("Null/Alternative Hypothesis") - rewrite the quoted hypothesis
# sentence in both paragraphs of the subtitle placeholder.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 3: Research question subtitle
# ---------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$rq = $s3.Shapes.Item(1).TextFrame.TextRange
$rqPara = $rq.Paragraphs(5, 1)

$rqPara.Text = '"Is there a correlation in the monthly average fund allocation to postpartum women over the period from October 2012 to September 2013 across all states?"'

# Re-split the new sentence into the same run layout the authors used,
# purely so the run boundaries line up with the source text.
$rqPara.Characters(1, 23)    # `"Is there a correlation`
$rqPara.Characters(24, 8)    # ` in the `
$rqPara.Characters(32, 8)    # `monthly `
$rqPara.Characters(40, 115)  # `average fund allocation to postpartum women over the period from October 2012 to September 2013 across all states?"`

# ---------------------------------------------------------------------
# Slide 4: Null / Alternative hypothesis subtitle
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$hyp = $s4.Shapes.Item(1).TextFrame.TextRange

# -- Null hypothesis (paragraph 1): keep "Null Hypothesis (H0): " run intact,
#    replace the quoted sentence that follows it.
$nullPara = $hyp.Paragraphs(1, 1)
$nullPrefixLen = $nullPara.Length - 133  # length of the existing quoted run before edit
$nullQuote = $nullPara.Characters($nullPrefixLen + 1, $nullPara.Length - $nullPrefixLen)
$nullQuote.Text = '"There is no correlation between time (October 2012 to September 2013) and the monthly average fund allocation to postpartum women."'

$nullPara = $hyp.Paragraphs(1, 1)
$nullPara.Characters(23, 79)    # `"There is no correlation between time (October 2012 to September 2013) and the `
$nullPara.Characters(102, 8)    # `monthly `
$nullPara.Characters(110, 45)   # `average fund allocation to postpartum women."`

# -- Alternative hypothesis (paragraph 2): keep "Alternative Hypothesis (H1): "
#    run intact, replace the quoted sentence that follows it.
$altPara = $hyp.Paragraphs(2, 1)
$altPrefixLen = $altPara.Length - 132
$altQuote = $altPara.Characters($altPrefixLen + 1, $altPara.Length - $altPrefixLen)
$altQuote.Text = '"There is a correlation between time (October 2012 to September 2013) and the monthly average fund allocation to postpartum women."'

$altPara = $hyp.Paragraphs(2, 1)
$altPara.Characters(30, 78)   # `"There is a correlation between time (October 2012 to September 2013) and the `
$altPara.Characters(108, 7)   # `monthly`
$altPara.Characters(115, 46)  # ` average fund allocation to postpartum women."`
